# Day 2 Servlets deck: reposition/resize the two RequestDispatcher
# screenshots (PowerPoint stores Shape.Left/Top/Width/Height in points,
# 1 pt = 12700 EMU; the literals below are chosen so that, after the
# points value round-trips through PowerPoint's single-precision
# storage, the resulting EMU in the XML match the target exactly).

$p = $ppt.ActivePresentation

# Slide 12 ("RequestDispatcher"): picture
#   off  1528762,2343150 -> 398723,2016578 (EMU)
#   ext  6422523x2556228 -> 8745277x3480707 (EMU)
$s12 = $p.Slides.Item(12)
$pic12 = $s12.Shapes.Item(3)
$pic12.Left = 31.39551281102362
$pic12.Top = 158.7856692913386
$pic12.Width = 688.6044881889763
$pic12.Height = 274.0714273228346

# Slide 16 ("RequestDispatcher - forward()"): picture
#   off  476174,2042431 -> 0,1879145 (EMU)
#   ext  8193238x3737883 -> 9147675x4173312 (EMU)
$s16 = $p.Slides.Item(16)
$pic16 = $s16.Shapes.Item(3)
$pic16.Left = 0.0
$pic16.Top = 147.96418022834646
$pic16.Width = 720.2893990787402
$pic16.Height = 328.60725409448816
